$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank columns, right-to-left so original column letters
# stay valid reference points while inserting.
# Final layout: new columns land at G, I, K, M.
$ws.Columns("J").Insert()
$ws.Columns("I").Insert()
$ws.Columns("H").Insert()
$ws.Columns("G").Insert()

# New headers (row 1) for the inserted columns.
$ws.Range("G1").Value = "Refinery-gasoline"
$ws.Range("I1").Value = "MtG-gasoline"
$ws.Range("K1").Value = "Gtkm-km"
$ws.Range("M1").Value = "B2gas-gasoline"

# New data values (row 2) mirroring the adjacent production column.
$ws.Range("G2").Value2 = $ws.Range("F2").Value2
$ws.Range("I2").Value2 = $ws.Range("H2").Value2
$ws.Range("K2").Value2 = $ws.Range("J2").Value2
$ws.Range("M2").Value2 = $ws.Range("L2").Value2
